$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Preis pro 5 Stück (Robot Arm)" label to "Preis pro 3 Stück (Robot Arm)"
$ws.Range("F24").Value = "Preis pro 3 Stück (Robot Arm)"

# Update the multiplier in the total formula from *5 to *3 to match the renamed label
$ws.Range("G24").Formula = "=SUM(G2:G17)*3 + G19"

# Move the active cell selection (cosmetic cursor position change)
$ws.Range("G25").Select()
